# Daily attendance processing - 2025-10-19 05:43:01
# Normalizes the "Recorded By" (column G) entries so that "System" is always
# listed first among the comma-separated recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ', '

    if ($parts.Length -ge 2 -and $parts[0] -ne 'System') {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newText = [string]::Join(', ', $parts)
        $cell.Value = $newText
    }
}
